$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D8").Value = -8.906300000000002
$ws.Range("D10").Value = -8.351099999999995
$ws.Range("D12").Value = -7.052799999999995
$ws.Range("D18").Value = -8.765099999999999
$ws.Range("D37").Value = -8.134499999999999
$ws.Range("D55").Value = -8.920899999999996
$ws.Range("D68").Value = -6.895999999999992
$ws.Range("D77").Value = -5.624100000000001
$ws.Range("D78").Value = -7.741900000000002
$ws.Range("D81").Value = -7.631499999999996
$ws.Range("D82").Value = -8.304399999999992
